# adding base functions, go to car brand, carrwall, fixing bugs
$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet (LoginTest) and name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "NewCarsTest"

# Header row
$ws.Range("A1").Value = "carBrand"
$ws.Range("B1").Value = "carTitle"

# Car brand rows
$ws.Range("A2").Value = "BMW"
$ws.Range("B2").Value = "BMW Cars"

$ws.Range("A3").Value = "Hyundai"
$ws.Range("B3").Value = "Hyundai Cars"

$ws.Range("A4").Value = "Toyota"
$ws.Range("B4").Value = "Toyota Cars"

$ws.Range("A5").Value = "Honda"
$ws.Range("B5").Value = "Honda Cars"

# Make the new sheet the active/selected sheet (like the author left it selected).
$ws.Activate()
